$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 4
$ws.Range("H4").Value = 1173.3914
$ws.Range("I4").Value = 762.3077
$ws.Range("K4").Value = 762.3077
$ws.Range("M4").Value = -648.3077

# Row 82
$ws.Range("H82").Value = 13827.728
$ws.Range("I82").Value = 11139.571
$ws.Range("J82").Value = 18532
$ws.Range("K82").Value = 33418.713
$ws.Range("L82").Value = 55596
$ws.Range("M82").Value = -33012.713
$ws.Range("N82").Value = -56408

# Row 85
$ws.Range("H85").Value = 13827.728
$ws.Range("I85").Value = 11139.571
$ws.Range("J85").Value = 18532
$ws.Range("K85").Value = 33418.713
$ws.Range("L85").Value = 55596
$ws.Range("M85").Value = -32014.713
$ws.Range("N85").Value = -58404

# Row 100
$ws.Range("H100").Value = 6411.8
$ws.Range("J100").Value = 6677.125
$ws.Range("L100").Value = 6677.125
$ws.Range("N100").Value = -7759.125

# Row 106
$ws.Range("H106").Value = 7344.5654
$ws.Range("I106").Value = 7344.5654
$ws.Range("K106").Value = 7344.5654
$ws.Range("M106").Value = -6713.5654

# Row 115
$ws.Range("H115").Value = 3041.5557
$ws.Range("I115").Value = 692
$ws.Range("J115").Value = 3712.8572
$ws.Range("K115").Value = 2076
$ws.Range("L115").Value = 11138.5716
$ws.Range("M115").Value = -509
$ws.Range("N115").Value = -14272.5716

# Row 137
$ws.Range("H137").Value = 2743.3704
$ws.Range("I137").Value = 2107.2354
$ws.Range("J137").Value = 3824.8
$ws.Range("K137").Value = 6321.706200000001
$ws.Range("L137").Value = 11474.4
$ws.Range("M137").Value = -3771.706200000001
$ws.Range("N137").Value = -16574.4

# Row 138
$ws.Range("H138").Value = 4425.0786
$ws.Range("I138").Value = 3152.889
$ws.Range("K138").Value = 9458.667000000001
$ws.Range("M138").Value = -4318.667000000001

# Row 141
$ws.Range("H141").Value = 7744.909
$ws.Range("I141").Value = 7744.909
$ws.Range("K141").Value = 23234.727
$ws.Range("M141").Value = -18054.727

$ws = $wb.Worksheets.Item("ARM")
# Row 45
$ws.Range("H45").Value = 1833659
$ws.Range("I45").Value = 2646824.2
$ws.Range("K45").Value = 2646824.2
$ws.Range("M45").Value = -2646447.2

# Row 61
$ws.Range("H61").Value = 6369036
$ws.Range("I61").Value = 7413257.5
$ws.Range("J61").Value = 1670038.4
$ws.Range("K61").Value = 7413257.5
$ws.Range("L61").Value = 1670038.4
$ws.Range("M61").Value = -7413045.5
$ws.Range("N61").Value = -1670462.4

# Row 63
$ws.Range("H63").Value = 1903.2084
$ws.Range("I63").Value = 2005.579
$ws.Range("J63").Value = 1514.2
$ws.Range("K63").Value = 2005.579
$ws.Range("L63").Value = 1514.2
$ws.Range("M63").Value = -1319.579
$ws.Range("N63").Value = -2886.2

# Row 66
$ws.Range("H66").Value = 1903.2084
$ws.Range("I66").Value = 2005.579
$ws.Range("J66").Value = 1514.2
$ws.Range("K66").Value = 10027.895
$ws.Range("L66").Value = 7571
$ws.Range("M66").Value = -6595.895
$ws.Range("N66").Value = -14435

# Row 102
$ws.Range("H102").Value = 2501.8076
$ws.Range("I102").Value = 1665.6316
$ws.Range("K102").Value = 1665.6316
$ws.Range("M102").Value = -43.63159999999993

# Row 136
$ws.Range("H136").Value = 6369036
$ws.Range("I136").Value = 7413257.5
$ws.Range("J136").Value = 1670038.4
$ws.Range("K136").Value = 22239772.5
$ws.Range("L136").Value = 5010115.199999999
$ws.Range("M136").Value = -22237222.5
$ws.Range("N136").Value = -5015215.199999999

$ws = $wb.Worksheets.Item("BSM")
# Row 94
$ws.Range("H94").Value = 1738.9714
$ws.Range("I94").Value = 2068.8096
$ws.Range("J94").Value = 1244.2142
$ws.Range("K94").Value = 2068.8096
$ws.Range("L94").Value = 1244.2142
$ws.Range("M94").Value = -1617.8096
$ws.Range("N94").Value = -2146.2142

# Row 105
$ws.Range("H105").Value = 648684.6
$ws.Range("I105").Value = 920906.8
$ws.Range("K105").Value = 920906.8
$ws.Range("M105").Value = -919159.8

$ws = $wb.Worksheets.Item("CRP")
# Row 31
$ws.Range("H31").Value = 25002386
$ws.Range("I31").Value = 27029028
$ws.Range("K31").Value = 27029028
$ws.Range("M31").Value = -27028733

# Row 34
$ws.Range("H34").Value = 25002386
$ws.Range("I34").Value = 27029028
$ws.Range("K34").Value = 27029028
$ws.Range("M34").Value = -27028826

# Row 132
$ws.Range("H132").Value = 2331.963
$ws.Range("I132").Value = 2058.6956
$ws.Range("J132").Value = 3903.25
$ws.Range("K132").Value = 6176.0868
$ws.Range("L132").Value = 11709.75
$ws.Range("M132").Value = -3646.0868
$ws.Range("N132").Value = -16769.75

# Row 141
$ws.Range("H141").Value = 381666.56
$ws.Range("J141").Value = 440714.16
$ws.Range("L141").Value = 440714.16
$ws.Range("N141").Value = -451074.16

$ws = $wb.Worksheets.Item("CUL")
# Row 62
$ws.Range("H62").Value = 13281.77
$ws.Range("I62").Value = 9256.5
$ws.Range("J62").Value = 16732
$ws.Range("K62").Value = 27769.5
$ws.Range("L62").Value = 50196
$ws.Range("M62").Value = -27083.5
$ws.Range("N62").Value = -51568

# Row 65
$ws.Range("H65").Value = 13281.77
$ws.Range("I65").Value = 9256.5
$ws.Range("J65").Value = 16732
$ws.Range("K65").Value = 83308.5
$ws.Range("L65").Value = 150588
$ws.Range("M65").Value = -79876.5
$ws.Range("N65").Value = -157452

# Row 69
$ws.Range("H69").Value = 17449.363
$ws.Range("I69").Value = 3633.3333
$ws.Range("J69").Value = 22630.375
$ws.Range("K69").Value = 10899.9999
$ws.Range("L69").Value = 67891.125
$ws.Range("M69").Value = -10088.9999
$ws.Range("N69").Value = -69513.125

# Row 72
$ws.Range("H72").Value = 17449.363
$ws.Range("I72").Value = 3633.3333
$ws.Range("J72").Value = 22630.375
$ws.Range("K72").Value = 32699.9997
$ws.Range("L72").Value = 203673.375
$ws.Range("M72").Value = -28643.9997
$ws.Range("N72").Value = -211785.375

# Row 107
$ws.Range("H107").Value = 4142151.2
$ws.Range("I107").Value = 6166
$ws.Range("K107").Value = 18498
$ws.Range("M107").Value = -16578

$ws = $wb.Worksheets.Item("GSM")
# Row 46
$ws.Range("H46").Value = 10000
$ws.Range("I46").Value = 10000
$ws.Range("K46").Value = 10000
$ws.Range("M46").Value = -9844

# Row 97
$ws.Range("H97").Value = 8335.823
$ws.Range("J97").Value = 14804.625
$ws.Range("L97").Value = 14804.625
$ws.Range("N97").Value = -15796.625

$ws = $wb.Worksheets.Item("LTW")
# Row 46
$ws.Range("H46").Value = 1090.4445
$ws.Range("I46").Value = 916.75
$ws.Range("J46").Value = 1229.4
$ws.Range("K46").Value = 916.75
$ws.Range("L46").Value = 1229.4
$ws.Range("M46").Value = -728.75
$ws.Range("N46").Value = -1605.4

# Row 68
$ws.Range("H68").Value = 1739402.4
$ws.Range("I68").Value = 2780747.5
$ws.Range("J68").Value = 3827.4443
$ws.Range("K68").Value = 2780747.5
$ws.Range("L68").Value = 3827.4443
$ws.Range("M68").Value = -2779998.5
$ws.Range("N68").Value = -5325.4443

# Row 71
$ws.Range("H71").Value = 1739402.4
$ws.Range("I71").Value = 2780747.5
$ws.Range("J71").Value = 3827.4443
$ws.Range("K71").Value = 13903737.5
$ws.Range("L71").Value = 19137.2215
$ws.Range("M71").Value = -13899993.5
$ws.Range("N71").Value = -26625.2215

$ws = $wb.Worksheets.Item("WVR")
# Row 63
$ws.Range("H63").Value = 34624.5
$ws.Range("J63").Value = 34624.5
$ws.Range("L63").Value = 34624.5
$ws.Range("N63").Value = -35872.5

# Row 66
$ws.Range("H66").Value = 34624.5
$ws.Range("J66").Value = 34624.5
$ws.Range("L66").Value = 103873.5
$ws.Range("N66").Value = -110113.5

# Row 74
$ws.Range("H74").Value = 20859
$ws.Range("I74").Value = 22084.5
$ws.Range("J74").Value = 20042
$ws.Range("K74").Value = 22084.5
$ws.Range("L74").Value = 20042
$ws.Range("M74").Value = -21148.5
$ws.Range("N74").Value = -21914

# Row 77
$ws.Range("H77").Value = 20859
$ws.Range("I77").Value = 22084.5
$ws.Range("J77").Value = 20042
$ws.Range("K77").Value = 66253.5
$ws.Range("L77").Value = 60126
$ws.Range("M77").Value = -61573.5
$ws.Range("N77").Value = -69486

# Row 126
$ws.Range("H126").Value = 3892.9412
$ws.Range("J126").Value = 910.5714
$ws.Range("L126").Value = 2731.7142
$ws.Range("N126").Value = -7671.7142

# Row 132
$ws.Range("H132").Value = 456602.38
$ws.Range("I132").Value = 2255.3684
$ws.Range("K132").Value = 6766.1052
$ws.Range("M132").Value = -4236.1052

# Row 136
$ws.Range("H136").Value = 327381.03
$ws.Range("I136").Value = 5201.304
$ws.Range("J136").Value = 1253647.8
$ws.Range("K136").Value = 15603.912
$ws.Range("L136").Value = 3760943.4
$ws.Range("M136").Value = -13053.912
$ws.Range("N136").Value = -3766043.4
